# Refresh cryptos list figures (price + 1h volume change) per upstream feed snapshot
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.874.38'
$ws.Range("E2").Value = '  +1.36%  '
$ws.Range("D3").Value = '2.352.33'
$ws.Range("E3").Value = '  +0.79%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.674'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.89%  '
$ws.Range("E6").Value = '  +2.69%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '72.50'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +10.43%  '
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.537'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +19.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0995'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '28.95'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +9.04%  '
$ws.Range("E12").Value = '  +2.65%  '
$ws.Range("D13").Value = '2.700.98'
$ws.Range("E13").Value = '  +0.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.74'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +9.49%  '
$ws.Range("E15").Value = '  +7.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.899'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +7.04%  '
$ws.Range("D17").Value = '2.355.47'
$ws.Range("E17").Value = '  +0.81%  '
$ws.Range("D18").Value = '43.881.03'
$ws.Range("E18").Value = '  +1.51%  '
$ws.Range("E19").Value = '  +4.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '77.85'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.51%  '
$ws.Range("E21").Value = '  +3.92%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '253.41'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.77'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.51'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.49'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.20'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.43'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '172.59'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.24%  '
$ws.Range("E30").Value = '  +6.56%  '
$ws.Range("E31").Value = '  +1.83%  '
$ws.Range("E32").Value = '  +5.90%  '
$ws.Range("E33").Value = '  +3.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0716'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.25'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.60%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.95'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.18%  '
$ws.Range("E37").Value = '  -2.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.41'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0266'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.41%  '
$ws.Range("E40").Value = '  +8.86%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("E42").Value = '  -1.21%  '
$ws.Range("E43").Value = '  +4.47%  '
$ws.Range("E44").Value = '  +4.24%  '
$ws.Range("E45").Value = '  +0.65%  '
$ws.Range("E46").Value = '  +1.26%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '98.07'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("E48").Value = '  +11.45%  '
$ws.Range("E49").Value = '  +4.14%  '
$ws.Range("D50").Value = '1.435.87'
$ws.Range("E50").Value = '  +0.28%  '
$ws.Range("E51").Value = '  +1.53%  '
